$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.554.35"
$ws.Range("E2").Value = "  -4.25%  "

$ws.Range("D3").Value = "3.111.45"
$ws.Range("E3").Value = "  -4.09%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").Value = "'552.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.36%  "

$ws.Range("D6").Value = "'137.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -11.04%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "

$ws.Range("D8").Value = "3.107.40"
$ws.Range("E8").Value = "  -3.97%  "

$ws.Range("E9").Value = "  -3.12%  "

$ws.Range("E10").Value = "  -4.98%  "

$ws.Range("D11").Value = "'6.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -9.87%  "

$ws.Range("D12").Value = "'0.476"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.88%  "

$ws.Range("D13").Value = "'35.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.82%  "

$ws.Range("D14").Value = "'0.0000219"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.55%  "

$ws.Range("D15").Value = "3.619.00"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D16").Value = "63.674.29"
$ws.Range("E16").Value = "  -4.20%  "

$ws.Range("E17").Value = "  -3.19%  "

$ws.Range("D18").Value = "3.116.00"
$ws.Range("E18").Value = "  -3.91%  "

$ws.Range("D19").Value = "'6.79"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.84%  "

$ws.Range("D20").Value = "'493.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -12.82%  "

$ws.Range("D21").Value = "'13.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.99%  "

$ws.Range("E22").Value = "  -2.72%  "

$ws.Range("D23").Value = "'7.28"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -8.01%  "

$ws.Range("D24").Value = "'79.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.51%  "

$ws.Range("D25").Value = "'12.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -8.75%  "

$ws.Range("E26").Value = "  +0.16%  "

$ws.Range("D27").Value = "'8.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -9.60%  "

$ws.Range("E28").Value = "  -6.07%  "

$ws.Range("D29").Value = "'1.99"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -12.29%  "

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.07%  "

$ws.Range("D31").Value = "'26.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.88%  "

$ws.Range("D32").Value = "'1.13"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.58%  "

$ws.Range("E33").Value = "  -8.64%  "

$ws.Range("D34").Value = "'59.08"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.71%  "

$ws.Range("D35").Value = "'522.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -7.86%  "

$ws.Range("D36").Value = "'6.03"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.18%  "

$ws.Range("D37").Value = "'5.19"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -10.21%  "

$ws.Range("E38").Value = "  -11.24%  "

$ws.Range("D39").Value = "3.154.67"
$ws.Range("E39").Value = "  +0.22%  "

$ws.Range("D40").Value = "'0.0808"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.00%  "

$ws.Range("D41").Value = "'0.121"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.25%  "

$ws.Range("E42").Value = "  -10.92%  "

$ws.Range("D43").Value = "'8.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.23%  "

$ws.Range("D44").Value = "'0.260"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.82%  "

$ws.Range("E45").Value = "  +0.06%  "

$ws.Range("D46").Value = "'2.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.76%  "

$ws.Range("D47").Value = "'25.30"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.59%  "

$ws.Range("D48").Value = "'121.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "

$ws.Range("E49").Value = "  -3.95%  "

$ws.Range("E50").Value = "  -9.17%  "

$ws.Range("E51").Value = "  -9.30%  "
